{"js": "// Replace each lattice-multiplication exercise cell with a new \"A x B\"\n// problem, regenerating the derived digit-breakdown lines that follow it.\n// Mapping is keyed by the ORIGINAL \"A x B\" problem text so the script is\n// resilient to table iteration order.\nconst replacements = {\n  \"25 x 12\": \"82 x 78\",\n  \"97 x 36\": \"67 x 95\",\n  \"44 x 70\": \"70 x 87\",\n  \"31 x 82\": \"12 x 75\",\n  \"13 x 35\": \"57 x 62\",\n  \"24 x 31\": \"41 x 12\",\n  \"75 x 51\": \"75 x 31\",\n  \"47 x 52\": \"13 x 52\",\n  \"34 x 96\": \"64 x 25\",\n  \"18 x 34\": \"96 x 24\",\n  \"44 x 83\": \"37 x 70\",\n  \"11 x 48\": \"81 x 50\",\n  \"69 x 24\": \"90 x 40\",\n  \"43 x 20\": \"77 x 87\",\n  \"33 x 92\": \"60 x 14\",\n};\n\nfunction buildCellText(problem) {\n  const parts = problem.split(\" x \");\n  const a = parts[0];\n  const b = parts[1];\n  const line1 = problem;\n  const line2 = \"  \" + b[0] + \"    \" + b[1];\n  const line3 = \"  ----\";\n  const line4 = a[0] + \"|    |\";\n  const line5 = a[1] + \"|    |\";\n  return [line1, line2, line3, line4, line5].join(\"\\v\");\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (let ti = 0; ti < tables.items.length; ti++) {\n  const table = tables.items[ti];\n  const paragraphs = table.body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  for (let pi = 0; pi < paragraphs.items.length; pi++) {\n    const para = paragraphs.items[pi];\n    const text = para.text;\n    const firstLine = text.split(\"\\v\")[0];\n    const newProblem = replacements[firstLine];\n    if (newProblem) {\n      para.insertText(buildCellText(newProblem), Word.InsertLocation.replace);\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each lattice-multiplication exercise cell with a new \"A x B\"\n# problem, regenerating the derived digit-breakdown lines that follow it.\n# Mapping is keyed by the ORIGINAL \"A x B\" problem text so the script is\n# resilient to table iteration order.\n# NOTE: string pieces here are built with the \"-f\" format operator rather\n# than \"+\" concatenation, because this PowerShell host coerces two\n# numeric-looking string operands of \"+\" into numeric addition (e.g.\n# \"  7    \" + \"8\" -> 15 instead of \"  7    8\").\n$replacements = @{\n    \"25 x 12\" = \"82 x 78\"\n    \"97 x 36\" = \"67 x 95\"\n    \"44 x 70\" = \"70 x 87\"\n    \"31 x 82\" = \"12 x 75\"\n    \"13 x 35\" = \"57 x 62\"\n    \"24 x 31\" = \"41 x 12\"\n    \"75 x 51\" = \"75 x 31\"\n    \"47 x 52\" = \"13 x 52\"\n    \"34 x 96\" = \"64 x 25\"\n    \"18 x 34\" = \"96 x 24\"\n    \"44 x 83\" = \"37 x 70\"\n    \"11 x 48\" = \"81 x 50\"\n    \"69 x 24\" = \"90 x 40\"\n    \"43 x 20\" = \"77 x 87\"\n    \"33 x 92\" = \"60 x 14\"\n}\n\n$d = $word.ActiveDocument\n$vbr = [char]11\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellText = $cell.Range.Text\n        $firstLine = $cellText.Split($vbr)[0]\n        if ($replacements.ContainsKey($firstLine)) {\n            $problem = $replacements[$firstLine]\n            $nums = $problem -split \" x \"\n            $a = $nums[0]\n            $b = $nums[1]\n            $line1 = $problem\n            $line2 = \"{0}{1}{2}{3}\" -f \"  \", $b.Substring(0,1), \"    \", $b.Substring(1,1)\n            $line3 = \"  ----\"\n            $line4 = \"{0}{1}\" -f $a.Substring(0,1), \"|    |\"\n            $line5 = \"{0}{1}\" -f $a.Substring(1,1), \"|    |\"\n            $newText = \"{0}{1}{2}{3}{4}{5}{6}{7}{8}\" -f $line1, $vbr, $line2, $vbr, $line3, $vbr, $line4, $vbr, $line5\n            $cell.Range.Text = $newText\n        }\n    }\n}\n"}
